$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the dataset. In the source sheet it
# lands at row 253, which pushes every existing record from row 253 down
# through row 303 down by one row (303 -> 304), growing the used range
# from A1:R303 to A1:R304.
$ws.Rows("253").Insert()

# Populate the newly inserted row 253 with the new record's data. Columns
# A, B, C, E, F, H, I, O and R hold the same constant values used by every
# other record in this sheet.
$ws.Range("A253").Value = 8
$ws.Range("B253").Value = "Terminal La Palmera de La Serena"
$ws.Range("C253").Value = "Coquimbo"
$ws.Range("D253").Value = 45005
$ws.Range("E253").Value = 4
$ws.Range("F253").Value = 100112037
$ws.Range("G253").Value = "Cebollín"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 1100
$ws.Range("K253").Value = 1000
$ws.Range("L253").Value = 1200
$ws.Range("M253").Value = 1100
$ws.Range("N253").Value = '$/paquete 6 unidades'
$ws.Range("O253").Value = "Provincia del Elquí"
$ws.Range("P253").Value = 183
$ws.Range("Q253").Value = 6
$ws.Range("R253").Value = "Hortaliza"
